{"js": "// Apply strikethrough formatting to the review-comment paragraphs that have\n// been addressed, leaving the still-open comments (and the \"S1:\"/\"S2:\"/\"S3:\"\n// section headers, the first \"CITE zenodo\" note, and a couple of untouched\n// comments) unformatted.\n//\n// Matching is done on a short, ASCII-only prefix of each paragraph's text so\n// we don't have to worry about curly-quote/encoding mismatches between this\n// source file and the document content.\nconst addressedPrefixes = [\n  \"198: one needs to pass (to is \",\n  \"241: The probability? \",\n  \"290: the numbers in this parag\",\n  \"Discussion first paragraph: is\",\n  \"366: the opening bracket is mi\",\n  \"382: delete the \",\n  \"388: delete the \",\n  \"432: \",\n  \"463: does my code about visual\",\n  \"first line: punctuation missin\",\n  \"1000 or 1,000? \",\n  \"Not sure what this refers to: \",\n  \"Figure S2.3: not sure what the\",\n  \"Figure S2.6: \",\n  \"The line above that same equat\",\n  \"691: delete \",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text;\n  if (addressedPrefixes.some((prefix) => text.indexOf(prefix) === 0)) {\n    paragraph.font.strikeThrough = true;\n  }\n}\n\n// Word silently drops the transient \"_GoBack\" bookmark (the last-edit-\n// position marker it maintains internally) the next time the document is\n// saved after an edit; remove it explicitly here to match.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Apply strikethrough formatting to the review-comment paragraphs that have\n# been addressed, leaving the still-open comments (and the \"S1:\"/\"S2:\"/\"S3:\"\n# section headers, the first \"CITE zenodo\" note, and a couple of untouched\n# comments) unformatted.\n#\n# Matching is done on a short, ASCII-only prefix of each paragraph's text so\n# we don't have to worry about curly-quote/encoding mismatches between this\n# source file and the document content.\n$addressedPrefixes = @(\n    \"198: one needs to pass (to is \",\n    \"241: The probability? \",\n    \"290: the numbers in this parag\",\n    \"Discussion first paragraph: is\",\n    \"366: the opening bracket is mi\",\n    \"382: delete the \",\n    \"388: delete the \",\n    \"432: \",\n    \"463: does my code about visual\",\n    \"first line: punctuation missin\",\n    \"1000 or 1,000? \",\n    \"Not sure what this refers to: \",\n    \"Figure S2.3: not sure what the\",\n    \"Figure S2.6: \",\n    \"The line above that same equat\",\n    \"691: delete \"\n)\n\n$d = $word.ActiveDocument\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $para = $d.Paragraphs($i)\n    $text = $para.Range.Text\n    foreach ($prefix in $addressedPrefixes) {\n        if ($text.StartsWith($prefix)) {\n            $para.Range.Font.StrikeThrough = 1\n            break\n        }\n    }\n}\n\n# Word silently drops the transient \"_GoBack\" bookmark (the last-edit-\n# position marker it maintains internally) the next time the document is\n# saved after an edit; remove it explicitly here to match.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
